# Adds a new "2022" column (S) to the transport statistics table, mirroring
# the formatting of the existing "2021" column (R), and selects cell T3
# (the cell to the right of the newly added header) as the active cell,
# matching the published workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - paste only the source formatting (used so the new column
# inherits the same cell styles already used by column R, instead of
# creating brand-new style entries).
$xlPasteFormats = -4122

# Row numbers that hold data in columns B:R that need a matching value
# added in column S, together with the new value for 2022.
$newValues = @{
    3  = 2022
    4  = 10444.200000000001
    5  = 21.7
    6  = 7361.6
    7  = 143.1
    8  = 844.2
    10 = "2 756,0"
    11 = "1 013,8"
    12 = "1 451,1"
    13 = 273.39999999999998
    14 = "-"
    15 = 17.7
}

foreach ($row in 3..15) {
    $srcCell = $ws.Range("R" + $row)
    $dstCell = $ws.Range("S" + $row)

    # Copy column R's formatting down into the new column S cell first, so
    # the new cell reuses the existing style (same as the rest of the row).
    $srcCell.Copy()
    $dstCell.PasteSpecial($xlPasteFormats)

    if ($newValues.ContainsKey($row)) {
        $dstCell.Value = $newValues[$row]
    }
}

# Row 9 (Таксопарк/Taxi row) has no 2022 figure published yet, so S9 stays
# empty but keeps the inherited number formatting from R9.

# Mirror the workbook's saved selection state (cell to the right of the new
# header cell, T3) exactly as captured in the published file.
$ws.Range("T3").Select()
